# Remove the "% Ridership vs Population" summary textboxes that were
# overlaid on the jurisdiction ridership screenshots (Phoenix, Tempe,
# Mesa slides). These shapes duplicated/obscured info already present
# elsewhere, so they are deleted entirely.

$p = $ppt.ActivePresentation

# Slide 6 - Phoenix Jurisdiction Bus Ridership with Population Growth
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item("TextBox 1").Delete()

# Slide 7 - Tempe Jurisdiction Bus Ridership with Population Growth
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item("TextBox 3").Delete()

# Slide 8 - Mesa Jurisdiction Bus Ridership with Population Growth
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item("TextBox 1").Delete()
